# Update NATMI LR-pair output (Fbn1-Itgb1) with recomputed TPM-based statistics.
# The underlying per-cluster ligand/receptor expression values were recomputed
# with a new TPM normalization; this updates the dependent columns
# (expressing-cell counts/rates, average & total expression, derived
# specificities, and the resulting edge weights/specificities) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs (sending) -> ECs (target)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 21.08397466666667
$ws.Range("H2").Value = 63.251924
$ws.Range("I2").Value = 0.06331801375981215
$ws.Range("J2").Value = 0.06331801375981214
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.29841822580325195
$ws.Range("P2").Value = 0.29841822580325195
$ws.Range("Q2").Value = 3544.4233407137076
$ws.Range("R2").Value = 31899.810066423368
$ws.Range("S2").Value = 0.018895249327589035
$ws.Range("T2").Value = 0.01889524932758903
# Row 3: ECs (sending) -> FAPs (target)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 21.08397466666667
$ws.Range("H3").Value = 63.251924
$ws.Range("I3").Value = 0.06331801375981215
$ws.Range("J3").Value = 0.06331801375981214
$ws.Range("O3").Value = 0.2893586437755395
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 3436.8193433046977
$ws.Range("R3").Value = 30931.37408974227
$ws.Range("S3").Value = 0.018321614588100193
$ws.Range("T3").Value = 0.01832161458810018
# Row 4: ECs (sending) -> MuSCs (target)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 21.08397466666667
$ws.Range("H4").Value = 63.251924
$ws.Range("I4").Value = 0.06331801375981215
$ws.Range("J4").Value = 0.06331801375981214
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 3499.8033813505735
$ws.Range("R4").Value = 31498.23043215516
$ws.Range("S4").Value = 0.018657381224343916
$ws.Range("T4").Value = 0.018657381224343912
# Row 5: ECs (sending) -> Resolving-Mac (target)
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 21.08397466666667
$ws.Range("H5").Value = 63.251924
$ws.Range("I5").Value = 0.06331801375981215
$ws.Range("J5").Value = 0.06331801375981214
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.11756162548016566
$ws.Range("P5").Value = 0.11756162548016565
$ws.Range("Q5").Value = 1396.322788940064
$ws.Range("R5").Value = 12566.905100460574
$ws.Range("S5").Value = 0.007443768619779012
$ws.Range("T5").Value = 0.007443768619779009
# Row 6: FAPs (sending) -> ECs (target)
$ws.Range("I6").Value = 0.8174956765497907
$ws.Range("J6").Value = 0.8174956765497907
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.29841822580325195
$ws.Range("P6").Value = 0.29841822580325195
$ws.Range("Q6").Value = 45761.87067217659
$ws.Range("R6").Value = 411856.83604958944
$ws.Range("S6").Value = 0.24395560939781766
$ws.Range("T6").Value = 0.24395560939781766
# Row 7: FAPs (sending) -> FAPs (target)
$ws.Range("I7").Value = 0.8174956765497907
$ws.Range("J7").Value = 0.8174956765497907
$ws.Range("O7").Value = 0.2893586437755395
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.23654944025881455
$ws.Range("T7").Value = 0.23654944025881444
# Row 8: FAPs (sending) -> MuSCs (target)
$ws.Range("I8").Value = 0.8174956765497907
$ws.Range("J8").Value = 0.8174956765497907
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 45185.78463123478
$ws.Range("R8").Value = 406672.0616811131
$ws.Range("S8").Value = 0.24088450633495742
$ws.Range("T8").Value = 0.24088450633495742
# Row 9: FAPs (sending) -> Resolving-Mac (target)
$ws.Range("I9").Value = 0.8174956765497907
$ws.Range("J9").Value = 0.8174956765497907
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.11756162548016566
$ws.Range("P9").Value = 0.11756162548016565
$ws.Range("Q9").Value = 18027.852979667357
$ws.Range("R9").Value = 162250.6768170062
$ws.Range("S9").Value = 0.09610612055820114
$ws.Range("T9").Value = 0.09610612055820113
# Row 10: MuSCs (sending) -> ECs (target)
$ws.Range("G10").Value = 39.60693866666666
$ws.Range("H10").Value = 118.820816
$ws.Range("I10").Value = 0.11894496778374845
$ws.Range("J10").Value = 0.11894496778374845
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.29841822580325195
$ws.Range("P10").Value = 0.29841822580325195
$ws.Range("Q10").Value = 6658.315620455254
$ws.Range("R10").Value = 59924.84058409731
$ws.Range("S10").Value = 0.03549534625425117
$ws.Range("T10").Value = 0.03549534625425117
# Row 11: MuSCs (sending) -> FAPs (target)
$ws.Range("G11").Value = 39.60693866666666
$ws.Range("H11").Value = 118.820816
$ws.Range("I11").Value = 0.11894496778374845
$ws.Range("J11").Value = 0.11894496778374845
$ws.Range("O11").Value = 0.2893586437755395
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 6456.177978333879
$ws.Range("R11").Value = 58105.60180500491
$ws.Range("S11").Value = 0.03441775456183069
$ws.Range("T11").Value = 0.034417754561830674
# Row 12: MuSCs (sending) -> MuSCs (target)
$ws.Range("G12").Value = 39.60693866666666
$ws.Range("H12").Value = 118.820816
$ws.Range("I12").Value = 0.11894496778374845
$ws.Range("J12").Value = 0.11894496778374845
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 6574.495561773491
$ws.Range("R12").Value = 59170.46005596144
$ws.Range("S12").Value = 0.035048503212323195
$ws.Range("T12").Value = 0.035048503212323195
# Row 13: MuSCs (sending) -> Resolving-Mac (target)
$ws.Range("G13").Value = 39.60693866666666
$ws.Range("H13").Value = 118.820816
$ws.Range("I13").Value = 0.11894496778374845
$ws.Range("J13").Value = 0.11894496778374845
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.11756162548016566
$ws.Range("P13").Value = 0.11756162548016565
$ws.Range("Q13").Value = 2623.0382048339607
$ws.Range("R13").Value = 23607.34384350565
$ws.Range("S13").Value = 0.013983363755343406
$ws.Range("T13").Value = 0.013983363755343405
# Row 14: Resolving-Mac (sending) -> ECs (target)
$ws.Range("G14").Value = 0.08036333333333333
$ws.Range("H14").Value = 0.24109
$ws.Range("I14").Value = 0.00024134190664861212
$ws.Range("J14").Value = 0.00024134190664861212
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.29841822580325195
$ws.Range("P14").Value = 0.29841822580325195
$ws.Range("Q14").Value = 13.509866090597773
$ws.Range("R14").Value = 121.58879481538
$ws.Range("S14").Value = 0.00007202082359405288
$ws.Range("T14").Value = 0.00007202082359405288
# Row 15: Resolving-Mac (sending) -> FAPs (target)
$ws.Range("G15").Value = 0.08036333333333333
$ws.Range("H15").Value = 0.24109
$ws.Range("I15").Value = 0.00024134190664861212
$ws.Range("J15").Value = 0.00024134190664861212
$ws.Range("O15").Value = 0.2893586437755395
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 13.099724452292223
$ws.Range("R15").Value = 117.89752007063
$ws.Range("S15").Value = 0.00006983436679404526
$ws.Range("T15").Value = 0.00006983436679404523
# Row 16: Resolving-Mac (sending) -> MuSCs (target)
$ws.Range("G16").Value = 0.08036333333333333
$ws.Range("H16").Value = 0.24109
$ws.Range("I16").Value = 0.00024134190664861212
$ws.Range("J16").Value = 0.00024134190664861212
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 13.339793382566665
$ws.Range("R16").Value = 120.0581404431
$ws.Range("S16").Value = 0.00007111416941842075
$ws.Range("T16").Value = 0.00007111416941842075
# Row 17: Resolving-Mac (sending) -> Resolving-Mac (target)
$ws.Range("G17").Value = 0.08036333333333333
$ws.Range("H17").Value = 0.24109
$ws.Range("I17").Value = 0.00024134190664861212
$ws.Range("J17").Value = 0.00024134190664861212
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.11756162548016566
$ws.Range("P17").Value = 0.11756162548016565
$ws.Range("Q17").Value = 5.322201126807778
$ws.Range("R17").Value = 47.89981014127
$ws.Range("S17").Value = 0.00002837254684209324
$ws.Range("T17").Value = 0.000028372546842093236
